# Adds the "resume" notification script (resume_1_* / resume_2_*) keys and
# their ru-RU / en-US translations, removes the stale "async.miku" state by
# simply not touching it (it was already absent), and refreshes the
# "count of keys" counter + active-sheet bookkeeping on the "main" sheet.

$wb = $excel.ActiveWorkbook

$keys = $wb.Worksheets.Item("keys")
$ruRU = $wb.Worksheets.Item("ru-RU")
$enUS = $wb.Worksheets.Item("en-US")
$main = $wb.Worksheets.Item("main")

# ---------------------------------------------------------------------
# 1) keys sheet: new key names resume_1_1 .. resume_2_2 (rows 5-9)
# ---------------------------------------------------------------------
$keys.Range("A5").Value = "resume_1_1"
$keys.Range("B5").Value = 5
$keys.Range("A6").Value = "resume_1_2"
$keys.Range("B6").Value = 5
$keys.Range("A7").Value = "resume_1_3"
$keys.Range("B7").Value = 5
$keys.Range("A8").Value = "resume_2_1"
$keys.Range("B8").Value = 5
$keys.Range("A9").Value = "resume_2_2"
$keys.Range("B9").Value = 5

# ---------------------------------------------------------------------
# 2) ru-RU sheet: translations for the keys above (rows 5-9), then the
#    last line of the second resume popup ("Тебе это понравилось?")
#    before the final key (resume_2_3) gets introduced.
# ---------------------------------------------------------------------
$ruRU.Range("B5").Value = "О, ты вернулся!"
$ruRU.Range("A5").Value = "resume_1_1"
$ruRU.Range("B6").Value = "Как дела?"
$ruRU.Range("A6").Value = "resume_1_2"
$ruRU.Range("B7").Value = "Т-ты скучал по мне?"
$ruRU.Range("A7").Value = "resume_1_3"
$ruRU.Range("B8").Value = "С возращением, Семпай❤!"
$ruRU.Range("A8").Value = "resume_2_1"
$ruRU.Range("B9").Value = "И?"
$ruRU.Range("A9").Value = "resume_2_2"
$ruRU.Range("B10").Value = "Тебе это понравилось?"

# ---------------------------------------------------------------------
# 3) keys sheet: the final key of the resume script.
# ---------------------------------------------------------------------
$keys.Range("A10").Value = "resume_2_3"
$keys.Range("B10").Value = 5

# ru-RU row 10 key (reuses the shared string created just above).
$ruRU.Range("A10").Value = "resume_2_3"

# ---------------------------------------------------------------------
# 4) en-US sheet: translations for all six new keys (rows 5-10).
# ---------------------------------------------------------------------
$enUS.Range("B5").Value = "Oh, you're back!"
$enUS.Range("A5").Value = "resume_1_1"
$enUS.Range("B6").Value = "How is it going?"
$enUS.Range("A6").Value = "resume_1_2"
$enUS.Range("B7").Value = "D-did you missed me?"
$enUS.Range("A7").Value = "resume_1_3"
$enUS.Range("B8").Value = "Welcome back, Sempai❤!"
$enUS.Range("A8").Value = "resume_2_1"
$enUS.Range("B9").Value = "And?"
$enUS.Range("A9").Value = "resume_2_2"
$enUS.Range("B10").Value = "Did you like it?"
$enUS.Range("A10").Value = "resume_2_3"

# ---------------------------------------------------------------------
# 5) Apply the center/center alignment used for the new rows (matches
#    the rest of each sheet's already-centered data rows).
# ---------------------------------------------------------------------
$keys.Range("A6:B9").HorizontalAlignment = -4108
$keys.Range("A6:B9").VerticalAlignment = -4108
$keys.Range("A10").HorizontalAlignment = -4108
$keys.Range("A10").VerticalAlignment = -4108

$ruRU.Range("A5:B9").HorizontalAlignment = -4108
$ruRU.Range("A5:B9").VerticalAlignment = -4108
$ruRU.Range("A10").HorizontalAlignment = -4108
$ruRU.Range("A10").VerticalAlignment = -4108

$enUS.Range("A5:B10").HorizontalAlignment = -4108
$enUS.Range("A5:B10").VerticalAlignment = -4108

# ---------------------------------------------------------------------
# 6) "main" sheet: bump the displayed "count of keys" (B2) 4 -> 10, and
#    leave the cursor where the author left it.
# ---------------------------------------------------------------------
$main.Range("B2").Value = 10
[void]$main.Range("D9").Select()

# ---------------------------------------------------------------------
# 7) Cursor / active-sheet bookkeeping to match the saved workbook view:
#    keys ends on B10, en-US ends on B10, ru-RU ends on B9 and is the
#    sheet left active when the file was saved.
# ---------------------------------------------------------------------
[void]$keys.Range("B10").Select()
[void]$enUS.Range("B10").Select()
[void]$ruRU.Range("B9").Select()
$ruRU.Activate()
